$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text so numeric-looking strings
# (e.g. "295.80", "0.999", thousand-dot formatted prices) are not coerced
# into numbers and do not lose formatting (trailing zeros, sci notation, etc).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.605.63"
$ws.Range("E2").Value = "  +6.08%  "

$ws.Range("D3").Value = "3.406.96"
$ws.Range("E3").Value = "  +6.83%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "578.28"
$ws.Range("E5").Value = "  +7.75%  "

$ws.Range("D6").Value = "155.69"
$ws.Range("E6").Value = "  +7.48%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.416.11"
$ws.Range("E8").Value = "  +6.92%  "

$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").Value = "7.54"
$ws.Range("E10").Value = "  +3.03%  "

$ws.Range("E11").Value = "  +7.73%  "

$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("D13").Value = "3.990.66"
$ws.Range("E13").Value = "  +6.66%  "

$ws.Range("E14").Value = "  +0.23%  "

$ws.Range("E15").Value = "  +7.89%  "

$ws.Range("D16").Value = "27.15"
$ws.Range("E16").Value = "  +5.63%  "

$ws.Range("D17").Value = "63.671.84"
$ws.Range("E17").Value = "  +6.17%  "

$ws.Range("D18").Value = "3.397.28"
$ws.Range("E18").Value = "  +7.06%  "

$ws.Range("D19").Value = "6.40"
$ws.Range("E19").Value = "  +2.61%  "

$ws.Range("D20").Value = "14.16"
$ws.Range("E20").Value = "  +7.08%  "

$ws.Range("D21").Value = "8.44"
$ws.Range("E21").Value = "  +3.13%  "

$ws.Range("D22").Value = "389.16"
$ws.Range("E22").Value = "  +5.48%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").Value = "0.537"
$ws.Range("E24").Value = "  +2.80%  "

$ws.Range("D25").Value = "71.17"
$ws.Range("E25").Value = "  +2.32%  "

$ws.Range("D26").Value = "9.58"
$ws.Range("E26").Value = "  +11.70%  "

$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").Value = "  +21.35%  "

$ws.Range("D28").Value = "0.182"
$ws.Range("E28").Value = "  +7.01%  "

$ws.Range("E29").Value = "  +1.44%  "

$ws.Range("E30").Value = "  +8.39%  "

$ws.Range("D31").Value = "6.60"
$ws.Range("E31").Value = "  +8.54%  "

$ws.Range("D32").Value = "5.77"
$ws.Range("E32").Value = "  +9.48%  "

$ws.Range("E33").Value = "  +13.38%  "

$ws.Range("D34").Value = "23.29"
$ws.Range("E34").Value = "  +3.70%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "6.73"
$ws.Range("E36").Value = "  +2.53%  "

$ws.Range("E37").Value = "  +10.07%  "

$ws.Range("D38").Value = "158.13"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  +11.80%  "

$ws.Range("D40").Value = "27.72"
$ws.Range("E40").Value = "  +5.39%  "

$ws.Range("D41").Value = "0.0773"
$ws.Range("E41").Value = "  +8.94%  "

$ws.Range("D42").Value = "2.929.37"
$ws.Range("E42").Value = "  +5.23%  "

$ws.Range("D43").Value = "0.0321"
$ws.Range("E43").Value = "  +4.23%  "

$ws.Range("D44").Value = "0.765"
$ws.Range("E44").Value = "  +6.53%  "

$ws.Range("D45").Value = "41.47"
$ws.Range("E45").Value = "  +3.89%  "

$ws.Range("D46").Value = "4.34"
$ws.Range("E46").Value = "  +3.26%  "

$ws.Range("E47").Value = "  +9.36%  "

$ws.Range("D48").Value = "22.65"
$ws.Range("E48").Value = "  +10.35%  "

$ws.Range("D49").Value = "3.451.09"
$ws.Range("E49").Value = "  +6.79%  "

$ws.Range("D50").Value = "6.36"
$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "295.80"
$ws.Range("E51").Value = "  +12.82%  "

# Restore default style (clears the quotePrefix/number-format bookkeeping
# Excel created while forcing text, without altering the actual General look).
$ws.Range("D2:E51").Style = "Normal"
